## Mock dados quando não consegue entrar na api
## - Renomeia a planilha existente para "V1"
## - Adiciona uma nova planilha "V1.1" com dados mock de pedidos
## - Marca todos os requisitos da V1 como "Concluido"
## - Acrescenta uma célula auxiliar A23 = 1

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename the original sheet, add the new "V1.1" sheet right after it
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "V1"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "V1.1"

# ---------------------------------------------------------------------------
# 2) V1: mark every requirement's Status (column H) as "Concluido"
# ---------------------------------------------------------------------------
for ($r = 4; $r -le 21; $r++) {
    $ws1.Range("H$r").Value = "Concluido"
}

# Row-height tweaks that came along with the edit
$ws1.Rows.Item(4).RowHeight = 36.75
$ws1.Rows.Item(5).RowHeight = 30

# New helper row at the bottom of V1
$ws1.Range("A23").Value = 1

# Keep the cursor where the author left it
$ws1.Range("F17").Select()

# ---------------------------------------------------------------------------
# 3) V1.1: mock data used when the API can't be reached
# ---------------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 3.3072916666666665
$ws2.Columns.Item(2).ColumnWidth = 22.877604166666668
$ws2.Columns.Item(3).ColumnWidth = 32.307291666666664
$ws2.Columns.Item(4).ColumnWidth = 9.451822916666666
$ws2.Columns.Item(5).ColumnWidth = 14.166666666666666
$ws2.Columns.Item(6).ColumnWidth = 11.307291666666666
$ws2.Columns.Item(7).ColumnWidth = 22.592447916666668
$ws2.Columns.Item(8).ColumnWidth = 5.592447916666667

# Header row (reuse V1's header / divider-row styles)
$headers = @("ID", "Nome do requisito", "Descrição", "Prioridade", "Data de criação", "Responsável")
$headerStyle = $ws1.Range("A2").Style
$dividerStyle = $ws1.Range("A3").Style
for ($c = 1; $c -le 6; $c++) {
    $cell = $ws2.Cells.Item(1, $c)
    $cell.Value = $headers[$c - 1]
    $cell.Style = $headerStyle
}
for ($c = 1; $c -le 6; $c++) {
    $ws2.Cells.Item(2, $c).Style = $dividerStyle
}
$ws2.Range("A1:A2").Merge()
$ws2.Range("B1:B2").Merge()
$ws2.Range("C1:C2").Merge()
$ws2.Range("D1:D2").Merge()
$ws2.Range("E1:E2").Merge()
$ws2.Range("F1:F2").Merge()

# Mock order-related requirements (RF1-RF6 of the order flow)
$rows = @(
    @("RF1", "Carregar pizzas ", "Carregar pizzas para novo pedido", "ALTA", 45591, "Guilherme"),
    @("RF2", "Cadastrar pedido", "Cadastrar pedido no sistema", "ALTA", 45591, "Guilherme"),
    @("RF3", "Consultar pedido", "Consultar a lista de pedidos", "ALTA", 45591, "Guilherme"),
    @("RF4", "Cancelar pedido", "Cancelar um pedido em andamento", "ALTA", 45591, "Guilherme"),
    @("RF5", "Carregar ultimos pedidos", "Carregar os ultimos pedidos feitos", "ALTA", 45591, "Guilherme"),
    @("RF6 ", "Login", "Entrar com as credenciais", "ALTA", 45591, "Guilherme")
)

# Establish the date style once (mm-dd-yy maps to Excel's built-in numFmtId 14)
$ws2.Cells.Item(3, 5).NumberFormat = "mm-dd-yy"
$dateStyle = $ws2.Cells.Item(3, 5).Style

$r = 3
foreach ($row in $rows) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    if ($r -eq 8) {
        # last row: the priority was typed with a leading apostrophe (quote-prefixed text)
        $ws2.Cells.Item($r, 4).Value = "'" + $row[3]
    } else {
        $ws2.Cells.Item($r, 4).Value = $row[3]
    }
    $ws2.Cells.Item($r, 5).Style = $dateStyle
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $ws2.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# Stray underlined/empty formatting left on H4 from the original author's edits
$ws2.Range("H4").Font.Underline = $true

$ws2.Range("H8").Select()
